$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "userID"
$ws.Range("B1").Value = "fName"
$ws.Range("C1").Value = "lName"
$ws.Range("D1").Value = "company"

# --- Row 2: John Smith / SAP ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "John"
$ws.Range("C2").Value = "Smith"
$ws.Range("D2").Value = "SAP"

# --- Row 3: Hans Miller / Deutsche Bank (company typed before the name) ---
$ws.Range("A3").Value = 2
$ws.Range("D3").Value = "Deutsche Bank"
$ws.Range("B3").Value = "Hans"
$ws.Range("C3").Value = "Miller"

# --- Rows 4 & 5: names filled first, companies filled in afterwards ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Guntram"
$ws.Range("C4").Value = "Schmitt"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Heinrich "
$ws.Range("C5").Value = "Vogel"

$ws.Range("D4").Value = "Sparkasse"
$ws.Range("D5").Value = "Commerzbank"

# --- Widen column D so "Deutsche Bank" / "Commerzbank" aren't clipped ---
$ws.Columns.Item(4).ColumnWidth = 16.83

# --- Leave the selection where the author left it when saving ---
$ws.Range("E8").Select()
